$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in L1/M1 from placeholder "{i}" to actual sample names
$ws.Range("L1").Value = "germline_test01"
$ws.Range("M1").Value = "germline_test02"

# Update the per-variant sample Counts column (N) to reflect the new sample set:
# rows 2-200 previously counted 2 samples, now count 1 (only germline_test01 carries it)
$ws.Range("N2:N200").Value = 1
# rows 201-244 previously counted 1 sample, now count 0 (no sample carries it)
$ws.Range("N201:N244").Value = 0
